$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("J2").Value = 3.2
$ws.Range("K2").Value = 3.35
$ws.Range("Q2").Value = 2.24

# Row 3 updates
$ws.Range("F3").Value = 1.84
$ws.Range("G3").Value = 1.93
$ws.Range("I3").Value = 6.4
$ws.Range("J3").Value = 3.2
$ws.Range("K3").Value = 3.55
$ws.Range("P3").Value = 1.51
$ws.Range("Q3").Value = 2.6
